# "Small tweaks to files"
#
# 1. Add a new "meta" worksheet (after "parameters") holding background /
#    provenance notes about the workbook, written by Lizzie Wolkovich (EMW).
# 2. Expand a few header/label cells on the "parameters" sheet so they also
#    name the Stan model / R script they refer to.
# 3. Update the view state on "parameters" (freeze pane moved, new selection)
#    and leave the new "meta" sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- New "meta" sheet, inserted right after "parameters" -------------------
$meta = $wb.Worksheets.Add($null, $ws1)
$meta.Name = "meta"
$meta.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$meta.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$meta.PageSetup.TopMargin = $excel.InchesToPoints(1)
$meta.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$meta.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$meta.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

$meta.Range("A1").Value = "Below metadata by Lizzie Wolkovich (EMW):"
$meta.Range("A1").Font.Bold = $true

$meta.Range("A3").Value = "On 3 January 2017:"
$meta.Range("A3").Font.Bold = $true

$meta.Range("A4").Value = "This file was started by Harold Eyster a couple weeks back to trouble-shoot Stan model building with (initially) fake data then real data."

$meta.Range("A5").Value = "It goes through the values given to the model from fake data and reports -- for different Stan models -- what the model returned. "

# --- Expand a few labels on "parameters" to reference the actual files -----
$ws1.Range("O1").Value = "WITH SPECIES SANS INTERACTIONS (germdate_sp_no-inter.stan)"
$ws1.Range("B2").Value = "Fake data assignment (germdate_fakedata.R)"
$ws1.Range("I1").Value = "WITH SPECIES (germdate_sp.stan)"
$ws1.Range("C1").Value = "WITHOUT SPECIES (germdate.stan)"

# --- Last line of the meta sheet, in red -----------------------------------
$meta.Range("A6").Value = "Red means …"
$meta.Range("A6").Font.Color = 255

# --- View state: "parameters" keeps a frozen pane but scrolled/selected ----
$ws1.Activate()
$ws1.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws1.Range("B7:C8").Select()

# --- "meta" ends up the active/selected sheet -------------------------------
$meta.Activate()
$meta.Range("A6").Select()
